$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 1812.2222
$ws.Range("I40").Value = 1627.826
$ws.Range("K40").Value = 1627.826
$ws.Range("M40").Value = -1452.826

$ws.Range("H69").Value = 10667.917
$ws.Range("I69").Value = 2013
$ws.Range("J69").Value = 12398.9
$ws.Range("K69").Value = 6039
$ws.Range("L69").Value = 37196.7
$ws.Range("M69").Value = -5165
$ws.Range("N69").Value = -38944.7

$ws.Range("H72").Value = 10667.917
$ws.Range("I72").Value = 2013
$ws.Range("J72").Value = 12398.9
$ws.Range("K72").Value = 18117
$ws.Range("L72").Value = 111590.1
$ws.Range("M72").Value = -13749
$ws.Range("N72").Value = -120326.1

$ws.Range("H132").Value = 2613.843
$ws.Range("I132").Value = 1348.7805
$ws.Range("J132").Value = 7800.6
$ws.Range("K132").Value = 4046.3415
$ws.Range("L132").Value = 23401.8
$ws.Range("M132").Value = -1516.3415
$ws.Range("N132").Value = -28461.8

$ws.Range("H137").Value = 3649.9167
$ws.Range("I137").Value = 3224.875
$ws.Range("J137").Value = 4500
$ws.Range("K137").Value = 9674.625
$ws.Range("L137").Value = 13500
$ws.Range("M137").Value = -7124.625
$ws.Range("N137").Value = -18600

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1395.9706
$ws.Range("I2").Value = 1314.9
$ws.Range("J2").Value = 1511.7858
$ws.Range("K2").Value = 1314.9
$ws.Range("L2").Value = 1511.7858
$ws.Range("M2").Value = -1201.9
$ws.Range("N2").Value = -1737.7858

$ws.Range("H32").Value = 3964.7048
$ws.Range("I32").Value = 3378.5532
$ws.Range("J32").Value = 5932.5
$ws.Range("K32").Value = 3378.5532
$ws.Range("L32").Value = 5932.5
$ws.Range("M32").Value = -3091.5532
$ws.Range("N32").Value = -6506.5

$ws.Range("H46").Value = 71322
$ws.Range("I46").Value = 5303.7144
$ws.Range("J46").Value = 129088
$ws.Range("K46").Value = 5303.7144
$ws.Range("L46").Value = 129088
$ws.Range("M46").Value = -4984.7144
$ws.Range("N46").Value = -129726

$ws.Range("H74").Value = 34492.3
$ws.Range("I74").Value = 59580.53
$ws.Range("J74").Value = 1684.6154
$ws.Range("K74").Value = 59580.53
$ws.Range("L74").Value = 1684.6154
$ws.Range("M74").Value = -58706.53
$ws.Range("N74").Value = -3432.6154

$ws.Range("H77").Value = 34492.3
$ws.Range("I77").Value = 59580.53
$ws.Range("J77").Value = 1684.6154
$ws.Range("K77").Value = 297902.65
$ws.Range("L77").Value = 8423.076999999999
$ws.Range("M77").Value = -293534.65
$ws.Range("N77").Value = -17159.077

$ws.Range("H116").Value = 1395.9706
$ws.Range("I116").Value = 1314.9
$ws.Range("J116").Value = 1511.7858
$ws.Range("K116").Value = 1314.9
$ws.Range("L116").Value = 1511.7858
$ws.Range("M116").Value = 979.0999999999999
$ws.Range("N116").Value = -6099.7858

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1395.9706
$ws.Range("I3").Value = 1314.9
$ws.Range("J3").Value = 1511.7858
$ws.Range("K3").Value = 1314.9
$ws.Range("L3").Value = 1511.7858
$ws.Range("M3").Value = -1200.9
$ws.Range("N3").Value = -1739.7858

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 885
$ws.Range("I16").Value = 727.75
$ws.Range("J16").Value = 1199.5
$ws.Range("K16").Value = 727.75
$ws.Range("L16").Value = 1199.5
$ws.Range("M16").Value = -440.75
$ws.Range("N16").Value = -1773.5

$ws.Range("H31").Value = 71430150
$ws.Range("I31").Value = 100000960
$ws.Range("J31").Value = 45456700
$ws.Range("K31").Value = 100000960
$ws.Range("L31").Value = 45456700
$ws.Range("M31").Value = -100000665
$ws.Range("N31").Value = -45457290

$ws.Range("H34").Value = 71430150
$ws.Range("I34").Value = 100000960
$ws.Range("J34").Value = 45456700
$ws.Range("K34").Value = 100000960
$ws.Range("L34").Value = 45456700
$ws.Range("M34").Value = -100000758
$ws.Range("N34").Value = -45457104

$ws.Range("H99").Value = 2514.4707
$ws.Range("I99").Value = 2200.9744
$ws.Range("K99").Value = 2200.9744
$ws.Range("M99").Value = -702.9744000000001

$ws.Range("H105").Value = 1251.4286
$ws.Range("I105").Value = 1251.4286
$ws.Range("J105").Value = 0
$ws.Range("K105").Value = 1251.4286
$ws.Range("L105").Value = 0
$ws.Range("M105").Value = 495.5714
$ws.Range("N105").ClearContents()

$ws.Range("H113").Value = 885
$ws.Range("I113").Value = 727.75
$ws.Range("J113").Value = 1199.5
$ws.Range("K113").Value = 727.75
$ws.Range("L113").Value = 1199.5
$ws.Range("M113").Value = 1442.25
$ws.Range("N113").Value = -5539.5

$ws.Range("H126").Value = 2514.4707
$ws.Range("I126").Value = 2200.9744
$ws.Range("K126").Value = 6602.9232
$ws.Range("M126").Value = -4132.9232

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H92").Value = 2500724.5
$ws.Range("J92").Value = 3334133.2
$ws.Range("L92").Value = 10002399.6
$ws.Range("N92").Value = -10004895.6

$ws.Range("H107").Value = 666942.7
$ws.Range("I107").Value = 245.83333
$ws.Range("J107").Value = 3333730
$ws.Range("K107").Value = 737.49999
$ws.Range("L107").Value = 10001190
$ws.Range("M107").Value = 1182.50001
$ws.Range("N107").Value = -10005030

$ws.Range("H113").Value = 676.2632
$ws.Range("I113").Value = 713.5454999999999
$ws.Range("J113").Value = 625
$ws.Range("K113").Value = 2140.6365
$ws.Range("L113").Value = 1875
$ws.Range("M113").Value = 29.36350000000039
$ws.Range("N113").Value = -6215

$ws.Range("H121").Value = 1691913.1
$ws.Range("I121").Value = 553.1667
$ws.Range("J121").Value = 2288863.8
$ws.Range("K121").Value = 1659.5001
$ws.Range("L121").Value = 6866591.399999999
$ws.Range("M121").Value = -349.5001
$ws.Range("N121").Value = -6869211.399999999

$ws.Range("H129").Value = 2058.04
$ws.Range("I129").Value = 1495.1666
$ws.Range("J129").Value = 2374.6562
$ws.Range("K129").Value = 4485.4998
$ws.Range("L129").Value = 7123.9686
$ws.Range("M129").Value = 514.5002000000004
$ws.Range("N129").Value = -17123.9686

$ws.Range("H131").Value = 860.3684
$ws.Range("J131").Value = 885.163
$ws.Range("L131").Value = 2655.489
$ws.Range("N131").Value = -12735.489

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 35886530
$ws.Range("I122").Value = 62500924
$ws.Range("K122").Value = 187502772
$ws.Range("M122").Value = -187500322
